$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-15 with the new trade data (values refreshed for "version 2").
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = "XAU_TRY"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1524
$ws.Cells.Item(2, 7).Value = 101524
$ws.Cells.Item(2, 8).Value = 1.92832

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "XAU_TRY"
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 2486
$ws.Cells.Item(3, 7).Value = 104011
$ws.Cells.Item(3, 8).Value = 2.85508

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "XAU_TRY"
$ws.Cells.Item(4, 4).Value = 9
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(4, 6).Value = 169
$ws.Cells.Item(4, 7).Value = 104180
$ws.Cells.Item(4, 8).Value = 0.56411

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = "XAU_TRY"
$ws.Cells.Item(5, 4).Value = 17
$ws.Cells.Item(5, 5).Value = 19
$ws.Cells.Item(5, 6).Value = 2428
$ws.Cells.Item(5, 7).Value = 107435
$ws.Cells.Item(5, 8).Value = 2.71832

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = "EUR_TRY"
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1378
$ws.Cells.Item(6, 7).Value = 101378
$ws.Cells.Item(6, 8).Value = 1.78212

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = "EUR_TRY"
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).Value = 11
$ws.Cells.Item(7, 6).Value = 248
$ws.Cells.Item(7, 7).Value = 102770
$ws.Cells.Item(7, 8).Value = 0.64372

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = "EUR_TRY"
$ws.Cells.Item(8, 4).Value = 12
$ws.Cells.Item(8, 5).Value = 14
$ws.Cells.Item(8, 6).Value = 1449
$ws.Cells.Item(8, 7).Value = 104219
$ws.Cells.Item(8, 8).Value = 1.81362

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = "GBP_TRY"
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 4761
$ws.Cells.Item(9, 7).Value = 104761
$ws.Cells.Item(9, 8).Value = 5.17153

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "GBP_TRY"
$ws.Cells.Item(10, 4).Value = 6
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 1774
$ws.Cells.Item(10, 7).Value = 106535
$ws.Cells.Item(10, 8).Value = 2.09807

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "GBP_TRY"
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = 14
$ws.Cells.Item(11, 6).Value = 2967
$ws.Cells.Item(11, 7).Value = 109503
$ws.Cells.Item(11, 8).Value = 3.19174

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = "GBP_TRY"
$ws.Cells.Item(12, 4).Value = 15
$ws.Cells.Item(12, 5).Value = 19
$ws.Cells.Item(12, 6).Value = 5749
$ws.Cells.Item(12, 7).Value = 115252
$ws.Cells.Item(12, 8).Value = 5.66171

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = "USD_TRY"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 6568
$ws.Cells.Item(13, 7).Value = 106568
$ws.Cells.Item(13, 8).Value = 6.98257

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = "USD_TRY"
$ws.Cells.Item(14, 4).Value = 10
$ws.Cells.Item(14, 5).Value = 14
$ws.Cells.Item(14, 6).Value = 2266
$ws.Cells.Item(14, 7).Value = 108834
$ws.Cells.Item(14, 8).Value = 2.53153

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = "USD_TRY"
$ws.Cells.Item(15, 4).Value = 16
$ws.Cells.Item(15, 5).Value = 19
$ws.Cells.Item(15, 6).Value = 3076
$ws.Cells.Item(15, 7).Value = 111910
$ws.Cells.Item(15, 8).Value = 3.23281

# Remove the now-obsolete trailing rows (16-23); the used range shrinks to A1:H15.
$ws.Rows("16:23").Delete()
